$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 409.67
$ws.Range("I15").Value = 409.67
$ws.Range("K15").Value = 1229.01
$ws.Range("M15").Value = -1060.01
$ws.Range("H33").Value = 1335
$ws.Range("I33").Value = 1019.6
$ws.Range("K33").Value = 1019.6
$ws.Range("M33").Value = -790.6
$ws.Range("H51").Value = 8335879.5
$ws.Range("I51").Value = 62502250
$ws.Range("J51").Value = 2591.4614
$ws.Range("K51").Value = 62502250
$ws.Range("L51").Value = 2591.4614
$ws.Range("M51").Value = -62501766
$ws.Range("N51").Value = -3559.4614
$ws.Range("H137").Value = 11321.14
$ws.Range("I137").Value = 18183
$ws.Range("K137").Value = 54549
$ws.Range("M137").Value = -51999
$ws.Range("H138").Value = 9496.958000000001
$ws.Range("I138").Value = 7524
$ws.Range("J138").Value = 9891.549999999999
$ws.Range("K138").Value = 22572
$ws.Range("L138").Value = 29674.65
$ws.Range("M138").Value = -17432
$ws.Range("N138").Value = -39954.64999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15956.777
$ws.Range("I32").Value = 15956.777
$ws.Range("K32").Value = 15956.777
$ws.Range("M32").Value = -15669.777
$ws.Range("H61").Value = 3764.0469
$ws.Range("I61").Value = 3185.2778
$ws.Range("K61").Value = 3185.2778
$ws.Range("M61").Value = -2973.2778
$ws.Range("H74").Value = 159861.69
$ws.Range("I74").Value = 178198.36
$ws.Range("K74").Value = 178198.36
$ws.Range("M74").Value = -177324.36
$ws.Range("H77").Value = 159861.69
$ws.Range("I77").Value = 178198.36
$ws.Range("K77").Value = 890991.7999999999
$ws.Range("M77").Value = -886623.7999999999
$ws.Range("H122").Value = 2976.077
$ws.Range("I122").Value = 2920.3044
$ws.Range("J122").Value = 3403.6667
$ws.Range("K122").Value = 8760.913199999999
$ws.Range("L122").Value = 10211.0001
$ws.Range("M122").Value = -6310.913199999999
$ws.Range("N122").Value = -15111.0001
$ws.Range("H136").Value = 3764.0469
$ws.Range("I136").Value = 3185.2778
$ws.Range("K136").Value = 9555.8334
$ws.Range("M136").Value = -7005.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 9689.875
$ws.Range("I99").Value = 10505
$ws.Range("K99").Value = 10505
$ws.Range("M99").Value = -9007

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 23144.75
$ws.Range("J28").Value = 29998.334
$ws.Range("L28").Value = 29998.334
$ws.Range("N28").Value = -30488.334
$ws.Range("H31").Value = 5865.696
$ws.Range("I31").Value = 3877.5217
$ws.Range("J31").Value = 7853.8696
$ws.Range("K31").Value = 3877.5217
$ws.Range("L31").Value = 7853.8696
$ws.Range("M31").Value = -3582.5217
$ws.Range("N31").Value = -8443.8696
$ws.Range("H34").Value = 5865.696
$ws.Range("I34").Value = 3877.5217
$ws.Range("J34").Value = 7853.8696
$ws.Range("K34").Value = 3877.5217
$ws.Range("L34").Value = 7853.8696
$ws.Range("M34").Value = -3675.5217
$ws.Range("N34").Value = -8257.8696
$ws.Range("H58").Value = 4938.65
$ws.Range("I58").Value = 5163.8823
$ws.Range("K58").Value = 5163.8823
$ws.Range("M58").Value = -4960.8823
$ws.Range("H132").Value = 10417.657
$ws.Range("I132").Value = 11999.018
$ws.Range("J132").Value = 2774.4167
$ws.Range("K132").Value = 35997.054
$ws.Range("L132").Value = 8323.250100000001
$ws.Range("M132").Value = -33467.054
$ws.Range("N132").Value = -13383.2501
$ws.Range("H136").Value = 4938.65
$ws.Range("I136").Value = 5163.8823
$ws.Range("K136").Value = 15491.6469
$ws.Range("M136").Value = -12941.6469
$ws.Range("H141").Value = 221246.14
$ws.Range("J141").Value = 231692.28
$ws.Range("L141").Value = 231692.28
$ws.Range("N141").Value = -242052.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 313.63635
$ws.Range("I7").Value = 335.3
$ws.Range("J7").Value = 97
$ws.Range("K7").Value = 1005.9
$ws.Range("L7").Value = 291
$ws.Range("M7").Value = -893.9000000000001
$ws.Range("N7").Value = -515
$ws.Range("H8").Value = 1013.6
$ws.Range("I8").Value = 1013.6
$ws.Range("K8").Value = 3040.8
$ws.Range("M8").Value = -2901.8
$ws.Range("H122").Value = 2476.3635
$ws.Range("I122").Value = 2407
$ws.Range("J122").Value = 2597.75
$ws.Range("K122").Value = 21663
$ws.Range("L122").Value = 23379.75
$ws.Range("M122").Value = -19213
$ws.Range("N122").Value = -28279.75
$ws.Range("H132").Value = 2647.077
$ws.Range("I132").Value = 1182
$ws.Range("K132").Value = 10638
$ws.Range("M132").Value = -8108

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 26501.5
$ws.Range("I57").Value = 26501.5
$ws.Range("K57").Value = 26501.5
$ws.Range("M57").Value = -25681.5
$ws.Range("H80").Value = 9369.111000000001
$ws.Range("I80").Value = 3665
$ws.Range("J80").Value = 16499.25
$ws.Range("K80").Value = 3665
$ws.Range("L80").Value = 16499.25
$ws.Range("M80").Value = -2667
$ws.Range("N80").Value = -18495.25
$ws.Range("H83").Value = 9369.111000000001
$ws.Range("I83").Value = 3665
$ws.Range("J83").Value = 16499.25
$ws.Range("K83").Value = 18325
$ws.Range("L83").Value = 82496.25
$ws.Range("M83").Value = -13333
$ws.Range("N83").Value = -92480.25
$ws.Range("H102").Value = 15787.352
$ws.Range("I102").Value = 19059.139
$ws.Range("J102").Value = 3927.125
$ws.Range("K102").Value = 19059.139
$ws.Range("L102").Value = 3927.125
$ws.Range("M102").Value = -17437.139
$ws.Range("N102").Value = -7171.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3460.077
$ws.Range("I22").Value = 2372
$ws.Range("K22").Value = 2372
$ws.Range("M22").Value = -2077
$ws.Range("H27").Value = 3460.077
$ws.Range("I27").Value = 2372
$ws.Range("K27").Value = 2372
$ws.Range("M27").Value = -2265
$ws.Range("H43").Value = 35862
$ws.Range("J43").Value = 35862
$ws.Range("L43").Value = 35862
$ws.Range("N43").Value = -36248
$ws.Range("H46").Value = 2507.077
$ws.Range("I46").Value = 1819.5
$ws.Range("K46").Value = 1819.5
$ws.Range("M46").Value = -1631.5
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value = 8486.272000000001
$ws.Range("J82").Value = 6294
$ws.Range("L82").Value = 6294
$ws.Range("N82").Value = -7016
$ws.Range("H85").Value = 8486.272000000001
$ws.Range("J85").Value = 6294
$ws.Range("L85").Value = 6294
$ws.Range("N85").Value = -8790
$ws.Range("H136").Value = 3347.8572
$ws.Range("I136").Value = 2541
$ws.Range("J136").Value = 4961.5713
$ws.Range("K136").Value = 7623
$ws.Range("L136").Value = 14884.7139
$ws.Range("M136").Value = -5073
$ws.Range("N136").Value = -19984.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 19999
$ws.Range("J28").Value = 19999
$ws.Range("L28").Value = 19999
$ws.Range("N28").Value = -20695
$ws.Range("H55").Value = 2999.6667
$ws.Range("I55").Value = 2999.6667
$ws.Range("K55").Value = 2999.6667
$ws.Range("M55").Value = -2722.6667
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H132").Value = 4906036
$ws.Range("I132").Value = 5685950.5
$ws.Range("K132").Value = 17057851.5
$ws.Range("M132").Value = -17055321.5
